$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1) Style/format plumbing first (before the values change so the
#    "donor" cells still carry the formats we want to clone), reusing
#    existing style entries instead of synthesizing new ones.
# -----------------------------------------------------------------------

# J1 should end up with the style the old I1 header cell had (s=4).
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# I1 header switches to the plain header style used by the rest of row 1 (s=3).
$ws.Range("A1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# New column J (rows 2-28) gets the style the APPUNTI column used to have (s=2).
$ws.Range("I2").Copy()
$ws.Range("J2:J28").PasteSpecial(-4122)

# APPUNTI column (I, rows 2-28) loses its special style and becomes a plain
# bordered cell like column A (s=1).
$ws.Range("A2").Copy()
$ws.Range("I2:I28").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# -----------------------------------------------------------------------
# 2) Header row text
# -----------------------------------------------------------------------
$ws.Range("J1").Value = "EFFETTIVO"

# -----------------------------------------------------------------------
# 3) Value fixes on existing rows 19-23
# -----------------------------------------------------------------------
$ws.Range("I19").Value = 1.5

$ws.Range("B20").Value = "30L"
$ws.Range("I20").Value = 3

$ws.Range("B21").Value = 30
$ws.Range("I21").Value = 0

$ws.Range("B22").Value = 26
$ws.Range("I22").Value = 1.5

$ws.Range("B23").Value = "30L"
$ws.Range("H23").Value = "bianco"
$ws.Range("I23").Value = 1.5

# -----------------------------------------------------------------------
# 4) New EFFETTIVO column values for rows 2-28
# -----------------------------------------------------------------------
$effettivo = @{
  2=1; 3=1; 4=1; 5=1; 6=1; 7=1; 8=1; 9=1; 10=1; 11=1; 12=1; 13=1
  14=1; 15=1; 16=1; 17=1; 18=1; 19=1; 20=1; 21=1; 22=1; 23=1
  24=0; 25=0; 26=1; 27=1; 28=1
}
foreach ($r in $effettivo.Keys) {
  $ws.Range("J$r").Value = $effettivo[$r]
}

# -----------------------------------------------------------------------
# 5) Rewritten rows 24-28 (degree-project rows shuffled/extended)
# -----------------------------------------------------------------------
$ws.Range("A24").Value = "Calcolo Numerico"
$ws.Range("B24").Value = 21
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = "Terzo"
$ws.Range("E24").Value = "Secondo"
$ws.Range("F24").Value = "Dedé"
$ws.Range("G24").Value = "matematica"
$ws.Range("H24").Value = "blu"
$ws.Range("I24").Value = 2

$ws.Range("A25").Value = "Fondamenti di Ricerca Operativa"
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = "Terzo"
$ws.Range("E25").Value = "Secondo"
$ws.Range("F25").Value = "Malucelli"
$ws.Range("G25").Value = "informatica"
$ws.Range("H25").Value = "bianco"
$ws.Range("I25").Value = 0

$ws.Range("A26").Value = "Progetto di Algoritmi e Principi dell'Informatica"
$ws.Range("B26").Value = "30L"
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = "Prova Finale"
$ws.Range("E26").Value = "."
$ws.Range("F26").Value = "Barenghi"
$ws.Range("G26").Value = "informatica"
$ws.Range("H26").Value = "."
$ws.Range("I26").Value = 1

$ws.Range("A27").Value = "Progetto di Reti Logiche"
$ws.Range("B27").Value = "30L"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "Prova Finale"
$ws.Range("E27").Value = "."
$ws.Range("F27").Value = "Fornaciari"
$ws.Range("G27").Value = "informatica"
$ws.Range("H27").Value = "."
$ws.Range("I27").Value = 1

$ws.Range("A28").Value = "Progetto di Ingegneria del Software"
$ws.Range("B28").Value = 27
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = "Prova Finale"
$ws.Range("E28").Value = "."
$ws.Range("F28").Value = "Margara"
$ws.Range("G28").Value = "informatica"
$ws.Range("H28").Value = "."
$ws.Range("I28").Value = 1

# -----------------------------------------------------------------------
# 6) Cosmetic bits: column widths, dimension/view follow automatically;
#    best-effort nudge of the two resized columns and the scroll/selection
#    state to match the saved workbook.
# -----------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 10.85546875
$ws.Columns.Item(10).ColumnWidth = 15.7109375

$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("I20").Select()
